$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row 17 with the same shape/style as row 16 above it.
$row = 17

$ws.Rows.Item($row - 1).Copy()
$ws.Rows.Item($row).PasteSpecial()
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 1).Value = 42622.885833333334

$ws.Cells.Item($row, 2).Value = 34
$ws.Cells.Item($row, 3).Value = 60
$ws.Cells.Item($row, 4).Value = 33
$ws.Cells.Item($row, 5).Value = 60
$ws.Cells.Item($row, 6).Value = 20
$ws.Cells.Item($row, 7).Value = 14707
$ws.Cells.Item($row, 8).Value = 28395
$ws.Cells.Item($row, 9).Value = 3137
$ws.Cells.Item($row, 10).Value = 464
$ws.Cells.Item($row, 11).Value = 260
$ws.Cells.Item($row, 12).Value = 47
$ws.Cells.Item($row, 13).Value = 12
$ws.Cells.Item($row, 14).Value = "Bag"

$wb.Save()
